$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table (Room ID / Min Nights) after adding newly-scraped rows and
# re-saving to a different folder. Header stays in row 1; data rows are
# rewritten in full so the shared-string table + row layout match.
$ws.Range("A1").Value = "Room ID"
$ws.Range("B1").Value = "Min Nights"

$ws.Range("A2").Value = "837352260137971048"
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "43435162"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "830193102361409290"
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = "49525472"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "716883705085875481"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "969823904189469776"
$ws.Range("B7").Value = 4

# Column A keeps its existing "text" number format, now right-aligned.
$ws.Range("A1:A7").NumberFormat = "@"
$ws.Range("A1:A7").HorizontalAlignment = -4152

# Column B: General format, right aligned, except the row that kept its
# original "0" numeric format (now row 4 - the original 830193102361409290
# row). B3 inherited the old "0" format from the original B3 cell, so
# clear it (and the brand-new rows) back to plain General first.
$ws.Range("B1").ClearFormats()
$ws.Range("B2").ClearFormats()
$ws.Range("B3").ClearFormats()
$ws.Range("B5").ClearFormats()
$ws.Range("B6").ClearFormats()
$ws.Range("B7").ClearFormats()
$ws.Range("B4").NumberFormat = "0"

$ws.Range("B1:B7").HorizontalAlignment = -4152

# Selection / active cell as saved in the file.
$ws.Range("B3").Select()
